$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("callReason")

# Append four new rows of lookup data to the callReason table
$ws.Range("A37").Value = "متابعة طلب صيانة"
$ws.Range("B37").Value = 36

$ws.Range("A38").Value = "متابعه شكوى"
$ws.Range("B38").Value = 37

$ws.Range("A39").Value = "ابلاغ موعد تسليم"
$ws.Range("B39").Value = 38

$ws.Range("A40").Value = "ابلاغ موعد سحب"
$ws.Range("B40").Value = 39

# Scroll the view down to show the newly added rows and select the last cell
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("A40").Select()

# Configure the page setup as part of implementing the printable view of the table
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
